$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (these must stay text, matching the
# source data which stores prices/volumes as literal strings).
$textCells = @(
    "D5", "D6", "D15", "D20", "D21", "D22", "D23", "D26",
    "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34",
    "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42",
    "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50",
    "D51"
)
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated values row by row.
$ws.Range("D2").Value = '92.179.23'
$ws.Range("E2").Value = '  +1.81%  '

$ws.Range("D3").Value = '3.110.63'
$ws.Range("E3").Value = '  -3.05%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '237.86'
$ws.Range("E5").Value = '  -0.84%  '

$ws.Range("D6").Value = '615.94'
$ws.Range("E6").Value = '  -0.59%  '

$ws.Range("E7").Value = '  -1.40%  '

$ws.Range("E8").Value = '  +5.16%  '

$ws.Range("E9").Value = '  -0.11%  '

$ws.Range("D10").Value = '3.107.93'
$ws.Range("E10").Value = '  -3.11%  '

$ws.Range("E11").Value = '  +0.73%  '

$ws.Range("E12").Value = '  -1.24%  '

$ws.Range("E13").Value = '  -0.40%  '

$ws.Range("D14").Value = '92.227.72'
$ws.Range("E14").Value = '  -8.44%  '

$ws.Range("D15").Value = '34.23'
$ws.Range("E15").Value = '  -3.51%  '

$ws.Range("E16").Value = '  -2.30%  '

$ws.Range("D17").Value = '3.705.48'
$ws.Range("E17").Value = '  -2.88%  '

$ws.Range("D18").Value = '3.126.96'
$ws.Range("E18").Value = '  -0.85%  '

$ws.Range("D20").Value = '14.63'
$ws.Range("E20").Value = '  -4.08%  '

$ws.Range("D21").Value = '5.80'
$ws.Range("E21").Value = '  -3.54%  '

$ws.Range("D22").Value = '9.40'
$ws.Range("E22").Value = '  +1.74%  '

$ws.Range("D23").Value = '446.56'
$ws.Range("E23").Value = '  -1.27%  '

$ws.Range("E24").Value = '  -3.81%  '

$ws.Range("E25").Value = '  -4.08%  '

$ws.Range("D26").Value = '86.73'
$ws.Range("E26").Value = '  -2.64%  '

$ws.Range("D27").Value = '11.82'
$ws.Range("E27").Value = '  -1.81%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.10%  '

$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = '0.133'
$ws.Range("E29").Value = '  -7.64%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.232'
$ws.Range("E30").Value = '  -0.15%  '

$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value = '0.171'
$ws.Range("E31").Value = '  -0.38%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '9.12'
$ws.Range("E32").Value = '  -3.13%  '

$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D33").Value = '7.88'
$ws.Range("E33").Value = '  +2.75%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.159'
$ws.Range("E34").Value = '  -6.78%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '26.16'
$ws.Range("E35").Value = '  -3.48%  '

$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").Value = '1.90'
$ws.Range("E36").Value = '  -4.01%  '

$ws.Range("B37").Value = 'MantraDAO'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D37").Value = '3.86'
$ws.Range("E37").Value = '  +1.85%  '

$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '483.23'
$ws.Range("E38").Value = '  -5.51%  '

$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = '1.29'
$ws.Range("E39").Value = '  -4.83%  '

$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").Value = '23.86'
$ws.Range("E40").Value = '  +8.07%  '

$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").Value = '0.433'
$ws.Range("E41").Value = '  -5.38%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '3.29'
$ws.Range("E42").Value = '  -4.63%  '

$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("B44").Value = 'Binance-PegBSC-USD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D44").Value = '0.726'
$ws.Range("E44").Value = '  -27.31%  '

$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").Value = '161.88'
$ws.Range("E45").Value = '  +3.42%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '1.89'
$ws.Range("E46").Value = '  -2.64%  '

$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '0.692'
$ws.Range("E47").Value = '  -6.27%  '

$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D48").Value = '1.39'
$ws.Range("E48").Value = '  -0.36%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0337'
$ws.Range("E49").Value = '  +3.92%  '

$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '4.43'
$ws.Range("E50").Value = '  -2.27%  '

$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = '43.93'
$ws.Range("E51").Value = '  -0.24%  '
